# Update cryptocurrency price (column D) and 1-hour volume change (column E)
# figures with the latest scraped values from the GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '26.877.26'
    'E2' = '  -1.58%  '
    'D3' = '1.825.71'
    'E3' = '  -1.72%  '
    'E4' = '  +0.54%  '
    'D5' = '311.00'
    'E5' = '  -0.84%  '
    'E6' = '  +0.57%  '
    'D7' = '0.4570'
    'E7' = '  -0.88%  '
    'D8' = '0.3688'
    'E8' = '  -0.40%  '
    'D9' = '0.07158'
    'E9' = '  -2.21%  '
    'D10' = '0.8737'
    'E10' = '  -1.08%  '
    'D11' = '0.07765'
    'E11' = '  -0.43%  '
    'D12' = '19.59'
    'E12' = '  -1.30%  '
    'D13' = '1.795.32'
    'E13' = '  -3.26%  '
    'D14' = '5.314'
    'E14' = '  -1.65%  '
    'D15' = '6.383'
    'E15' = '  -2.58%  '
    'D16' = '86.71'
    'E16' = '  -5.60%  '
    'D17' = '1.009'
    'E17' = '  +0.74%  '
    'D18' = '0.000008712'
    'E18' = '  -1.65%  '
    'E19' = '  +0.49%  '
    'D20' = '26.918.07'
    'E20' = '  -1.50%  '
    'D21' = '14.46'
    'E21' = '  -2.25%  '
    'D22' = '5.003'
    'E22' = '  -2.32%  '
    'D23' = '2.067.88'
    'E23' = '  +0.88%  '
    'D24' = '10.42'
    'E24' = '  -0.59%  '
    'D25' = '2.001'
    'E25' = '  +4.75%  '
    'D26' = '151.42'
    'E26' = '  -0.37%  '
    'D27' = '18.16'
    'E27' = '  -1.04%  '
    'D28' = '1.962'
    'E28' = '  -4.96%  '
    'D29' = '113.66'
    'E29' = '  -2.06%  '
    'D30' = '4.910'
    'E30' = '  -4.18%  '
    'D31' = '0.08799'
    'E31' = '  -0.54%  '
    'D32' = '3.048'
    'E32' = '  +1.59%  '
    'D33' = '0.7487'
    'E33' = '  -2.55%  '
    'E34' = '  -0.19%  '
    'D35' = '1.132'
    'E35' = '  -3.47%  '
    'D36' = '2.544'
    'E36' = '  -3.12%  '
    'D37' = '1.084'
    'E37' = '  +0.98%  '
    'D38' = '0.01945'
    'E38' = '  -1.00%  '
    'E39' = '  -1.53%  '
    'E40' = '  -1.06%  '
    'D41' = '6.938'
    'E41' = '  -1.60%  '
    'D42' = '0.4967'
    'E42' = '  -3.55%  '
    'D43' = '0.1597'
    'E43' = '  -2.71%  '
    'D44' = '8.316'
    'E44' = '  -1.04%  '
    'D45' = '0.4686'
    'E45' = '  -3.12%  '
    'D46' = '1.007'
    'E46' = '  +0.70%  '
    'D47' = '10.14'
    'E47' = '  -1.38%  '
    'D48' = '101.85'
    'E48' = '  -1.20%  '
    'D49' = '1.611'
    'E49' = '  -2.55%  '
    'D50' = '0.06102'
    'E50' = '  -1.84%  '
    'D51' = '64.49'
    'E51' = '  -1.84%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text formatting so numeric-looking values (e.g. "311.00")
    # are kept verbatim instead of being normalised by Excel, then
    # restore the default "Normal" style so no extra formatting sticks.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
